# The commit removes the first data record (id=1, "Edwards-Thomas Incubator")
# from the Incubators sheet. Every row below it shifts up by one: the id
# numbering becomes contiguous (2..30 -> 1..29 effectively becomes 2..30 one
# row higher), the former row 31 ("Chavez PLC Incubator") becomes the new
# last row (30), and the sheet shrinks from 31 to 30 data rows.
#
# Deleting the entire row 2 reproduces exactly that: it removes the old row
# 2's values and shifts rows 3:31 up into 2:30, leaving no row 31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(2).Delete()
